# Generate Report for Handoff
# Swap the reporting rows for 4671043e-... and 5a546f92-... (the 5a546f92 entry is now
# listed first / "ready for handoff" status moves onto the 4671043e entry) across the
# Overview / zh-cn / de-de sheets, and widen the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$idA = "4671043e-79a3-44df-ba54-b798b1604ef4"
$idB = "5a546f92-ab30-4fcf-bbec-bfe5f7c73a43"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "$idB.md"
$ws1.Range("B2").Value = "e2e\$idB.md"

$ws1.Range("A3").Value = "$idA.md"
$ws1.Range("B3").Value = "e2e\$idA.md"

$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-09-04 18:53:28"

$ws1.Hyperlinks.Delete()
$null = $ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idA.md", "", "", "e2e\$idB.md")
$null = $ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idB.md", "", "", "e2e\$idA.md")

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "$idB.md"
$ws2.Range("G2").Value = "$idB.d9abb552313603295ad77449c1f87efd3a2edee5.zh-cn.xlf"
$ws2.Range("I2").Value = "$idB.md"
$ws2.Range("J2").Value = "$idB.d9abb552313603295ad77449c1f87efd3a2edee5.zh-cn.xlf"

$ws2.Range("A3").Value = "$idA.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "$idA.4ddb0be76cb358a7c0ce4470de5500fd3755ce1d.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-04 18:53:24"
$ws2.Range("I3").Value = "$idA.md"
$ws2.Range("J3").Value = "$idA.4ddb0be76cb358a7c0ce4470de5500fd3755ce1d.zh-cn.xlf"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idA.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/505a707a78660fabcaf1107f7f0e70fbb3834e90/e2e/$idA.md."

$ws2.Columns.Item(16).ColumnWidth = 39.1666666666667

$ws2.Hyperlinks.Delete()
$null = $ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idA.md", "", "", "$idB.md")
$null = $ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fddc803e26bee6c045cb1219f69b3293759f0393/e2e/$idA.md", "", "", "$idB.md")
$null = $ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idB.md", "", "", "$idA.md")
$null = $ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fddc803e26bee6c045cb1219f69b3293759f0393/e2e/$idB.md", "", "", "$idA.md")

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "$idB.md"
$ws3.Range("G2").Value = "$idB.d9abb552313603295ad77449c1f87efd3a2edee5.de-de.xlf"
$ws3.Range("I2").Value = "$idB.md"
$ws3.Range("J2").Value = "$idB.d9abb552313603295ad77449c1f87efd3a2edee5.de-de.xlf"

$ws3.Range("A3").Value = "$idA.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "$idA.4ddb0be76cb358a7c0ce4470de5500fd3755ce1d.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-04 18:53:28"
$ws3.Range("I3").Value = "$idA.md"
$ws3.Range("J3").Value = "$idA.4ddb0be76cb358a7c0ce4470de5500fd3755ce1d.de-de.xlf"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idA.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/505a707a78660fabcaf1107f7f0e70fbb3834e90/e2e/$idA.md."

$ws3.Columns.Item(16).ColumnWidth = 39.1666666666667

$ws3.Hyperlinks.Delete()
$null = $ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idA.md", "", "", "$idB.md")
$null = $ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/25942bf638ffa3212db89859201d279fc8c3f9df/e2e/$idA.md", "", "", "$idB.md")
$null = $ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dae4197392d746f947c8871cb467bf5b6d498dc3/e2e/$idB.md", "", "", "$idA.md")
$null = $ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/25942bf638ffa3212db89859201d279fc8c3f9df/e2e/$idB.md", "", "", "$idA.md")
